# Hindalco prices update (2026-02-10 05:48:58 UTC)
# Rows 3-14 get new Basic Price / Circular Date values; their Circular Link
# (hyperlink + URL text) is removed entirely (cell becomes blank).
# All other rows' hyperlinks must be preserved exactly as before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: capture every existing hyperlink's target URL, keyed by the
# worksheet row it lives on (hyperlinks only ever sit in column F here).
# ---------------------------------------------------------------------
$urls = @{}
$hlCount = $ws.Hyperlinks.Count()
for ($i = 1; $i -le $hlCount; $i++) {
    $row = $i + 1
    $cellRef = "F" + $row
    $urls[$row] = $ws.Range($cellRef).Value()
}

# ---------------------------------------------------------------------
# Step 2: wipe every hyperlink (this runtime only supports deleting the
# whole collection reliably) - we will re-create the ones we still need.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# Step 3: apply the new Basic Price (D) / Circular Date (E) values for
# rows 3-14 and blank out their Circular Link (F) cell.
# ---------------------------------------------------------------------
$newData = @{
    3  = @{ D = 326.5;  E = "07.02.2026" }
    4  = @{ D = 326.5;  E = "07.02.2026" }
    5  = @{ D = 326.5;  E = "07.02.2026" }
    6  = @{ D = 323.5;  E = "06.02.2026" }
    7  = @{ D = 327;    E = "05.02.2026" }
    8  = @{ D = 332.25; E = "04.02.2026" }
    9  = @{ D = 330.75; E = "03.02.2026" }
    10 = @{ D = 338.25; E = "31.01.2026" }
    11 = @{ D = 338.25; E = "31.01.2026" }
    12 = @{ D = 338.25; E = "31.01.2026" }
    13 = @{ D = 360;    E = "30.01.2026" }
    14 = @{ D = 352.25; E = "29.01.2026" }
}

foreach ($row in $newData.Keys) {
    $info = $newData[$row]

    $ws.Range("D" + $row).Value2 = $info.D

    # Force the circular-date cell to stay plain text (it looks like a
    # date, and this engine - like Excel - will otherwise silently turn
    # it into a date serial number).
    $eCell = $ws.Range("E" + $row)
    $eCell.NumberFormat = "@"
    $eCell.Value2 = $info.E

    $ws.Range("F" + $row).ClearContents()
}

# ---------------------------------------------------------------------
# Step 4: restore the hyperlinks for every row that still needs one
# (everything except rows 3-14). Re-applying Hyperlinks.Add switches the
# cell onto the built-in "Hyperlink" style, so reset it back to the
# original plain centered look used throughout this sheet.
# ---------------------------------------------------------------------
foreach ($row in $urls.Keys) {
    if ($row -ge 3 -and $row -le 14) {
        continue
    }

    $target = $ws.Range("F" + $row)
    $ws.Hyperlinks.Add($target, $urls[$row])

    $target.Style = "Normal"
    $target.HorizontalAlignment = -4108
    $target.VerticalAlignment = -4108
}
